$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Weekly/28-day/YTD crime statistics table (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -66.666666666666
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = -14.285714285714
$ws.Range("I14").Value = 31
$ws.Range("J14").Value = 24
$ws.Range("K14").Value = 29.166666666666
$ws.Range("L14").Value = -8.823529411764
$ws.Range("M14").Value = 14.814814814814
$ws.Range("N14").Value = -83.510638297872
# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -71.428571428571
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = -18.181818181818
$ws.Range("I15").Value = 73
$ws.Range("J15").Value = 103
$ws.Range("K15").Value = -29.126213592233
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -32.407407407407
$ws.Range("N15").Value = -72.962962962963
# Row 16
$ws.Range("C16").Value = 41
$ws.Range("D16").Value = 41
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 170
$ws.Range("G16").Value = 166
$ws.Range("H16").Value = 2.409638554216
$ws.Range("I16").Value = 1006
$ws.Range("J16").Value = 1091
$ws.Range("K16").Value = -7.791017415215
$ws.Range("L16").Value = 14.971428571428
$ws.Range("M16").Value = -17.877551020408
$ws.Range("N16").Value = -80.094974277799
# Row 17
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 63
$ws.Range("E17").Value = -4.761904761904
$ws.Range("F17").Value = 248
$ws.Range("G17").Value = 273
$ws.Range("H17").Value = -9.157509157509
$ws.Range("I17").Value = 1594
$ws.Range("J17").Value = 1656
$ws.Range("K17").Value = -3.743961352657
$ws.Range("L17").Value = 10.311418685121
$ws.Range("M17").Value = 55.512195121951
$ws.Range("N17").Value = -50
# Row 18
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 36
$ws.Range("E18").Value = -19.444444444444
$ws.Range("F18").Value = 116
$ws.Range("G18").Value = 125
$ws.Range("H18").Value = -7.2
$ws.Range("I18").Value = 832
$ws.Range("J18").Value = 985
$ws.Range("K18").Value = -15.532994923857
$ws.Range("L18").Value = 12.280701754386
$ws.Range("M18").Value = 16.853932584269
$ws.Range("N18").Value = -85.682326621923
# Row 19
$ws.Range("C19").Value = 119
$ws.Range("D19").Value = 138
$ws.Range("E19").Value = -13.768115942029
$ws.Range("F19").Value = 514
$ws.Range("G19").Value = 552
$ws.Range("H19").Value = -6.884057971014
$ws.Range("I19").Value = 3270
$ws.Range("J19").Value = 3432
$ws.Range("K19").Value = -4.720279720279
$ws.Range("L19").Value = 27.684498242873
$ws.Range("M19").Value = 33.360522022838
$ws.Range("N19").Value = -43.949262941378
# Row 20
$ws.Range("C20").Value = 26
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = -18.75
$ws.Range("F20").Value = 139
$ws.Range("G20").Value = 113
$ws.Range("H20").Value = 23.008849557522
$ws.Range("I20").Value = 732
$ws.Range("J20").Value = 695
$ws.Range("K20").Value = 5.323741007194
$ws.Range("L20").Value = 60.175054704595
$ws.Range("M20").Value = 144.816053511706
$ws.Range("N20").Value = -85.647058823529
# Row 21
$ws.Range("C21").Value = 278
$ws.Range("D21").Value = 320
$ws.Range("E21").Value = -13.125
$ws.Range("F21").Value = 1202
$ws.Range("G21").Value = 1247
$ws.Range("H21").Value = -3.608660785886
$ws.Range("I21").Value = 7538
$ws.Range("J21").Value = 7986
$ws.Range("K21").Value = -5.609817180065
$ws.Range("L21").Value = 21.855803427093
$ws.Range("M21").Value = 28.898768809849
$ws.Range("N21").Value = -70.375319316172
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = -75
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 26
$ws.Range("H22").Value = -26.923076923076
$ws.Range("I22").Value = 156
$ws.Range("J22").Value = 172
$ws.Range("K22").Value = -9.302325581395
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = 24.8
# Row 23
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 24
$ws.Range("E23").Value = -20.833333333333
$ws.Range("F23").Value = 94
$ws.Range("G23").Value = 102
$ws.Range("H23").Value = -7.843137254901
$ws.Range("I23").Value = 655
$ws.Range("J23").Value = 686
$ws.Range("K23").Value = -4.518950437317
$ws.Range("L23").Value = 1.708074534161
$ws.Range("M23").Value = 50.229357798165
# Row 24
$ws.Range("C24").Value = 265
$ws.Range("D24").Value = 295
$ws.Range("E24").Value = -10.169491525423
$ws.Range("F24").Value = 1198
$ws.Range("G24").Value = 1290
$ws.Range("H24").Value = -7.131782945736
$ws.Range("I24").Value = 7600
$ws.Range("J24").Value = 8472
$ws.Range("K24").Value = -10.292728989612
$ws.Range("L24").Value = 20.520139549635
$ws.Range("M24").Value = 59.865376525031
# Row 25
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = 73
$ws.Range("E25").Value = 16.438356164383
$ws.Range("F25").Value = 350
$ws.Range("G25").Value = 362
$ws.Range("H25").Value = -3.314917127071
$ws.Range("I25").Value = 2470
$ws.Range("J25").Value = 2505
$ws.Range("K25").Value = -1.397205588822
$ws.Range("L25").Value = 15.528531337698
$ws.Range("M25").Value = -16.185951815405
# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -21.052631578947
$ws.Range("I26").Value = 139
$ws.Range("J26").Value = 162
$ws.Range("K26").Value = -14.197530864197
$ws.Range("L26").Value = 6.923076923076
# Row 27
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = -18.75
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 61
$ws.Range("H27").Value = -9.83606557377
$ws.Range("I27").Value = 329
$ws.Range("J27").Value = 379
$ws.Range("K27").Value = -13.192612137203
$ws.Range("L27").Value = 1.543209876543
# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 25
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 83
$ws.Range("J28").Value = 108
$ws.Range("K28").Value = -23.148148148148
$ws.Range("L28").Value = -36.641221374045
$ws.Range("M28").Value = -22.429906542056
$ws.Range("N28").Value = -81.348314606741
# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 50
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 75
$ws.Range("J29").Value = 91
$ws.Range("K29").Value = -17.582417582417
$ws.Range("L29").Value = -35.897435897435
$ws.Range("M29").Value = -19.354838709677
$ws.Range("N29").Value = -81.751824817518
# Row 30
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 200
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -16.666666666666
$ws.Range("I30").Value = 36
$ws.Range("J30").Value = 61
$ws.Range("K30").Value = -40.983606557377
$ws.Range("L30").Value = -18.181818181818

# --- Row 30 (Hate Crimes): D30/E30 change from text placeholders to real numbers ---
# Apply the same number formats used by sibling cells in the table (C30/F30 etc. and H30/K30/L30 etc.)
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
